$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows with new Fecha (D), Volumen (J), Precio minimo/maximo/promedio (K/L/M),
# Origen (O) and Precio $/Kg (P) values per the weekly refresh of this
# "Hortaliza, Macroferia Regional de Talca - Cilantro" dataset.

$ws.Range("D2").Value = 44355
$ws.Range("J2").Value = 150
$ws.Range("O2").Value = 'Región Metropolitana'
$ws.Range("D3").Value = 44725
$ws.Range("D5").Value = 44715
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 9000
$ws.Range("P5").Value = 250
$ws.Range("D6").Value = 44719
$ws.Range("D7").Value = 44364
$ws.Range("J7").Value = 100
$ws.Range("D8").Value = 44376
$ws.Range("D9").Value = 44701
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 7000
$ws.Range("O9").Value = 'Región del Maule'
$ws.Range("P9").Value = 194
$ws.Range("D10").Value = 44358
$ws.Range("O10").Value = 'Región Metropolitana'
$ws.Range("D11").Value = 44729
$ws.Range("D12").Value = 44362
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 6500
$ws.Range("L12").Value = 6500
$ws.Range("M12").Value = 6500
$ws.Range("P12").Value = 181
$ws.Range("D13").Value = 44371
$ws.Range("J13").Value = 150
$ws.Range("D14").Value = 44340
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 7000
$ws.Range("M14").Value = 7000
$ws.Range("O14").Value = 'Región del Maule'
$ws.Range("P14").Value = 194
$ws.Range("D15").Value = 44342
$ws.Range("D16").Value = 44720
$ws.Range("D17").Value = 44372
$ws.Range("O17").Value = 'Región Metropolitana'
$ws.Range("D18").Value = 44706
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 9000
$ws.Range("L18").Value = 9000
$ws.Range("M18").Value = 9000
$ws.Range("P18").Value = 250
$ws.Range("D19").Value = 44726
$ws.Range("J19").Value = 150
$ws.Range("K19").Value = 8000
$ws.Range("L19").Value = 8000
$ws.Range("M19").Value = 8000
$ws.Range("P19").Value = 222
$ws.Range("D20").Value = 44711
$ws.Range("K20").Value = 8500
$ws.Range("L20").Value = 8500
$ws.Range("M20").Value = 8500
$ws.Range("O20").Value = 'Región Metropolitana'
$ws.Range("P20").Value = 236
$ws.Range("D21").Value = 44348
$ws.Range("K21").Value = 7000
$ws.Range("L21").Value = 7000
$ws.Range("M21").Value = 7000
$ws.Range("O21").Value = 'Región del Maule'
$ws.Range("P21").Value = 194
$ws.Range("D22").Value = 44386
$ws.Range("K22").Value = 6500
$ws.Range("L22").Value = 6500
$ws.Range("M22").Value = 6500
$ws.Range("P22").Value = 181
$ws.Range("D23").Value = 44690
$ws.Range("J23").Value = 500
$ws.Range("O23").Value = 'Región del Maule'
$ws.Range("D24").Value = 44707
$ws.Range("J24").Value = 150
$ws.Range("K24").Value = 9000
$ws.Range("L24").Value = 9000
$ws.Range("M24").Value = 9000
$ws.Range("P24").Value = 250
